$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row (row 1) strings: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the used range A1:U93 into an Excel Table ("Table1")
$rng = $ws.Range("A1:U93")
$lo = $ws.ListObjects.Add(1, $rng, $false, $true)
$lo.Name = "Table1"

# 3. Freeze the header row (freeze pane below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
